$d = $word.ActiveDocument

# 1. Update semester/year on the title page: "Fall 2021" -> "Spring 2022"
$d.Content.Find.Execute("Fall 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Spring 2022", 2) | Out-Null

# 2. Fix capitalization of "webcampus" -> "WebCampus" everywhere it occurs
#    (appears before example_1, example_2/example_3, and example_4 download
#    instructions). A single whole-document ReplaceAll covers every instance.
$d.Content.Find.Execute("webcampus", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "WebCampus", 2) | Out-Null

# 3. Remove stray period after "achieved?" before "Is there any difference..."
$d.Content.Find.Execute("achieved?.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "achieved?", 2) | Out-Null

# 4. Remove comma: "pressed, and" -> "pressed and"
$d.Content.Find.Execute("pressed, and", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "pressed and", 2) | Out-Null
